$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" column (H) with the same header formatting as the
# existing header cells (B1:G1) - copy formats from G1 onto H1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 1
